# feat: add 2022-Q1 data
#
# The old "总计" (totals) sheet is repurposed into the new "2022-Q1"
# fund-holding sheet (it keeps its original sheetId/rId), and a brand new
# "总计" sheet is appended at the end with the refreshed totals table
# (original rows shifted down by one, with a new 2022-Q1 row on top).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Turn the previous "总计" sheet into the new "2022-Q1" sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Header row (B1:H1) — reuse B1's existing bordered header style for the
# three brand-new header cells (E1:H1).
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data rows. Column A keeps its existing 0-based index / style from the
# old sheet. Numeric-looking text (fund codes, and the D/E/F/G columns
# which are stored as text, not numbers) is entered with a leading
# apostrophe so Excel keeps it as text instead of coercing to a number.
$q1Data = @(
    @("005823", "泰康颐享混合A", "14.39", "20.19", "1.31", "0.1885", 6),
    @("519198", "万家颐和灵活配置混合", "1.78", "91.03", "5.43", "0.0967", 5),
    @("005824", "泰康颐享混合C", "2.82", "20.19", "1.31", "0.0369", 6),
    @("011765", "兴银高端制造混合A", "1.01", "93.23", "2.73", "0.0276", 7),
    @("011766", "兴银高端制造混合C", "0.39", "93.23", "2.73", "0.0106", 7)
)

$r = 2
foreach ($row in $q1Data) {
    $q1.Range("B$r").Value = "'" + $row[0]
    $q1.Range("C$r").Value = $row[1]
    $q1.Range("D$r").Value = "'" + $row[2]
    $q1.Range("E$r").Value = "'" + $row[3]
    $q1.Range("F$r").Value = "'" + $row[4]
    $q1.Range("G$r").Value = "'" + $row[5]
    $q1.Range("H$r").Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Append a brand new "总计" sheet after "2022-Q1" with the refreshed
#    totals table.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$total = $wb.Worksheets.Add($null, $lastSheet)
$total.Name = "总计"

# Reuse header / index-column styling from the 2022-Q1 sheet.
$q1.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$q1.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalData = @(
    @("2022-Q1", 5, 0.36),
    @("2021-Q4", 17, 2.61),
    @("2021-Q3", 27, 2),
    @("2021-Q2", 2, 0.06),
    @("2021-Q1", 5, 0.59),
    @("2020-Q4", 3, 0.11)
)

$r = 2
$idx = 0
foreach ($row in $totalData) {
    $total.Range("A$r").Value = $idx
    $total.Range("B$r").Value = $row[0]
    $total.Range("C$r").Value = $row[1]
    $total.Range("D$r").Value = $row[2]
    $r = $r + 1
    $idx = $idx + 1
}

Write-Output "done"
